$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 3;   Q = 33; R = 23 },
    @{ Row = 10;  Q = 55; R = 37 },
    @{ Row = 17;  Q = 35; R = 18 },
    @{ Row = 23;  Q = 21; R = 11 },
    @{ Row = 32;  Q = 97; R = 32 },
    @{ Row = 40;  Q = 75; R = 1  },
    @{ Row = 49;  Q = 55; R = 16 },
    @{ Row = 58;  Q = 22; R = 4  },
    @{ Row = 66;  Q = 10; R = 2  },
    @{ Row = 74;  Q = 44; R = 16 },
    @{ Row = 78;  Q = 24; R = 10 },
    @{ Row = 89;  Q = 50; R = 20 },
    @{ Row = 97;  Q = 41; R = 13 },
    @{ Row = 106; Q = 30; R = 22 },
    @{ Row = 115; Q = 34; R = 19 },
    @{ Row = 124; Q = 60; R = 55 },
    @{ Row = 133; Q = 15; R = 2  },
    @{ Row = 142; Q = 83; R = 62 }
)

foreach ($u in $updates) {
    $ws.Range("Q$($u.Row)").Value = $u.Q
    $ws.Range("R$($u.Row)").Value = $u.R
}
